# "after contact test case"
# Update the ContactUs sheet: store the phone number as text with a leading
# "+91" country code, format the cell as Text, resize column A to fit, and
# move the active selection to B8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the phone-number column (A) as Text so the leading "+" is preserved,
# then write the new value. Setting the number format before the value keeps
# Excel from re-interpreting the string as a number.
$ws.Range("A1:A2").NumberFormat = "@"
$ws.Range("A2").Value = "+917788445511"

# Resize column A to fit the new, longer text.
$ws.Columns("A").ColumnWidth = 13.3

# Acknowledge/ignore the "number stored as text" warning Excel raises for the
# leading "+" phone number (mirrors clicking "Ignore Error" in the UI).
try {
    $ws.Range("A2").Errors.Item(9).Ignore = $true
} catch {
}

# Move/record the active selection as it was left in the authored workbook.
$null = $ws.Range("B8").Select()
